# Auto-update predictions and index for 2025-11-01
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New fixture data for rows 2-19 (columns A-G)
# Order: Fixture, Pick, Competition, Time, Win Info, Confidence %, Odds
$data = @(
    @("Brighton vs Leeds", "Brighton", "England", "1st Nov 15:00", "96 of 123 Tips", "78", "2.00"),
    @("Nottingham Forest vs Man Utd", "Man Utd", "England", "1st Nov 15:00", "95 of 140 Tips", "68", "2.05"),
    @("Fulham vs Wolves", "Fulham", "England", "1st Nov 15:00", "61 of 91 Tips", "67", "1.80"),
    @("Heidenheim vs Eintracht Frankfurt", "Eintracht Frankfurt", "Germany", "1st Nov 14:30", "60 of 65 Tips", "92", "1.85"),
    @("Liverpool vs Aston Villa", "Liverpool", "England", "1st Nov 20:00", "59 of 86 Tips", "69", "1.65"),
    @("Burnley vs Arsenal", "Arsenal", "England", "1st Nov 15:00", "59 of 71 Tips", "83", "1.26"),
    @("Napoli vs Como", "Napoli", "Italy", "1st Nov 17:00", "44 of 58 Tips", "76", "1.91"),
    @("RB Leipzig vs Stuttgart", "RB Leipzig", "Germany", "1st Nov 14:30", "44 of 57 Tips", "77", "2.10"),
    @("Real Madrid vs Valencia", "Real Madrid", "Spain", "1st Nov 20:00", "44 of 48 Tips", "92", "1.18"),
    @("Monaco vs Paris FC", "Monaco", "France", "1st Nov 18:00", "43 of 51 Tips", "84", "1.67"),
    @("Cremonese vs Juventus", "Juventus", "Italy", "1st Nov 19:45", "42 of 49 Tips", "86", "1.52"),
    @("Auxerre vs Marseille", "Marseille", "France", "1st Nov 20:05", "39 of 45 Tips", "87", "1.70"),
    @("Udinese vs Atalanta", "Atalanta", "Italy", "1st Nov 14:00", "33 of 53 Tips", "62", "1.77"),
    @("PSG vs Nice", "PSG", "France", "1st Nov 16:00", "28 of 31 Tips", "90", "1.20"),
    @("Atletico Madrid vs Sevilla", "Atletico Madrid", "Spain", "1st Nov 15:15", "27 of 33 Tips", "82", "1.35"),
    @("Ajax vs Heerenveen", "Ajax", "Netherlands", "1st Nov 15:30", "26 of 29 Tips", "90", "1.55"),
    @("Bayern Munich vs Bayer Leverkusen", "Bayern Munich", "Germany", "1st Nov 17:30", "25 of 30 Tips", "83", "1.23"),
    @("Huachipato vs Universidad de Chile", "Universidad de Chile", "Chile", "2nd Nov 15:30", "25 of 26 Tips", "96", "1.73")
)

# Columns A-E are plain text and assign cleanly via .Value.
# Columns F (Confidence %) and G (Odds) are numeric-looking strings that must
# stay stored as text (as in the source inlineStr cells), so they are entered
# with a leading apostrophe to force text interpretation, matching how Excel
# itself keeps "numbers stored as text" verbatim (e.g. "78", "2.00").
$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = "'" + $rec[5]
    $ws.Cells.Item($row, 7).Value = "'" + $rec[6]
    $row++
}

# Clear the old H19 average formula (the average row moved down to H20)
$ws.Cells.Item(19, 8).ClearContents()

# New average row moved from row 19 to row 20
$ws.Range("H20").Formula = "=AVERAGE(F2:F19)"
